$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.821.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.86%  "

$ws.Range("D3").Value = "'2.621.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.28%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'597.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.80%  "

$ws.Range("D6").Value = "'151.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.08%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  +0.66%  "

$ws.Range("D9").Value = "'0.109"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.27%  "

$ws.Range("E10").Value = "  +3.31%  "

$ws.Range("E11").Value = "  +2.92%  "

$ws.Range("E12").Value = "  -0.83%  "

$ws.Range("D13").Value = "'27.96"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.60%  "

$ws.Range("D14").Value = "'3.097.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.37%  "

$ws.Range("D15").Value = "'63.744.49"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.01%  "

$ws.Range("D16").Value = "'0.0000154"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.67%  "

$ws.Range("D17").Value = "'2.630.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.02%  "

$ws.Range("D18").Value = "'12.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.98%  "

$ws.Range("D19").Value = "'4.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.78%  "

$ws.Range("D20").Value = "'348.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.64%  "

$ws.Range("D21").Value = "'6.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.87%  "

$ws.Range("D22").Value = "'1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.11%  "

$ws.Range("D23").Value = "'67.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.58%  "

$ws.Range("E24").Value = "  +7.52%  "

$ws.Range("B25").Value = "Fetch.AI"
$ws.Range("C25").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D25").Value = "'1.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.53%  "

$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").Value = "'9.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.97%  "

$ws.Range("D27").Value = "'556.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.77%  "

$ws.Range("D28").Value = "'8.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.58%  "

$ws.Range("E29").Value = "  +0.46%  "

$ws.Range("D31").Value = "'2.06"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.72%  "

$ws.Range("D32").Value = "'0.0₃0856"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.24%  "

$ws.Range("D33").Value = "'1.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.63%  "

$ws.Range("D34").Value = "'5.30"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.22%  "

$ws.Range("D35").Value = "'166.98"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.75%  "

$ws.Range("D36").Value = "'0.417"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.73%  "

$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.03%  "

$ws.Range("D38").Value = "'19.61"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.32%  "

$ws.Range("E39").Value = "  +0.42%  "

$ws.Range("E40").Value = "  +0.01%  "

$ws.Range("D41").Value = "'168.01"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.33%  "

$ws.Range("D42").Value = "'39.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.44%  "

$ws.Range("D43").Value = "'3.99"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.59%  "

$ws.Range("E44").Value = "  +4.94%  "

$ws.Range("D45").Value = "'22.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.80%  "

$ws.Range("D46").Value = "'0.636"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.87%  "

$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").Value = "'2.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.46%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0252"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.98%  "

$ws.Range("D49").Value = "'0.0971"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.51%  "

$ws.Range("D50").Value = "'19.40"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.79%  "

$ws.Range("D51").Value = "'0.0₆0240"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +22.12%  "
